$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.459.70"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "'1.550.72"
$ws.Range("E3").Value = "  -2.13%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "'210.77"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("E6").Value = "  -1.66%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").Value = "'23.92"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("E10").Value = "  -1.50%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "'1.771.67"
$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").Value = "'1.548.97"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").Value = "'28.440.85"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("E15").Value = "  -2.27%  "

$ws.Range("D16").Value = "'0.509"
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("D17").Value = "'60.99"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("D18").Value = "'228.37"
$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D23").Value = "'8.90"
$ws.Range("E23").Value = "  -2.75%  "

$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("D25").Value = "'150.79"
$ws.Range("E25").Value = "  -0.79%  "

$ws.Range("D26").Value = "'14.74"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").Value = "'6.21"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").Value = "'0.0467"
$ws.Range("E30").Value = "  -3.39%  "

$ws.Range("E31").Value = "  -4.44%  "

$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").Value = "'1.382.43"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").Value = "'3.00"
$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("E36").Value = "  -2.94%  "

$ws.Range("D37").Value = "'2.29"
$ws.Range("E37").Value = "  -3.00%  "

$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("E40").Value = "  +1.36%  "

$ws.Range("D41").Value = "'0.511"
$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").Value = "'0.771"
$ws.Range("E43").Value = "  -2.53%  "

$ws.Range("E44").Value = "  -1.27%  "

$ws.Range("D45").Value = "'5.33"
$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("D46").Value = "'61.74"
$ws.Range("E46").Value = "  -2.24%  "

$ws.Range("D47").Value = "'1.685.15"
$ws.Range("E47").Value = "  -2.07%  "

$ws.Range("D48").Value = "'0.876"
$ws.Range("E48").Value = "  -8.90%  "

$ws.Range("D49").Value = "'85.11"
$ws.Range("E49").Value = "  -1.82%  "

$ws.Range("D50").Value = "'43.13"
$ws.Range("E50").Value = "  +8.52%  "

$ws.Range("D51").Value = "'0.0₆0100"
$ws.Range("E51").Value = "  -2.44%  "
